# evo: new merge fields
# Replace the legacy "[onshow.XXX]" merge-field placeholders in column B of
# Feuil1 with the new "[namespace.field]" syntax, and move the active
# selection from D13 to B29 (matching the author's final cursor position).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newValues = @{
    2  = "[contact.title]"
    3  = "[contact.lastname]"
    4  = "[contact.firstname]"
    5  = "[contact.society]"
    6  = "[contact.adrs_num]"
    7  = "[contact.adrs_street]"
    8  = "[contact.adrs_comp]"
    9  = "[contact.adrs_town]"
    10 = "[contact.adrs_pc]"
    11 = "[contact.adrs_contry]"
    12 = "[user.lastname]"
    13 = "[user.firstname]"
    14 = "[res_letterbox.destination]"
    15 = "[res_letterbox.doctype]"
    16 = "[res_letterbox.category_id]"
    17 = "[res_letterbox.nature]"
    18 = "[res_letterbox.admission_date]"
    19 = "[res_letterbox.doc_date]"
    20 = "[res_letterbox.process_limit_date]"
    21 = "[res_letterbox.process_notes]"
    22 = "[res_letterbox.closing_date]"
    23 = "[res_letterbox.subject]"
    24 = "[res_letterbox.chrono]"
    25 = "[res_letterbox.author]"
    26 = "[res_letterbox.creation_date]"
    27 = "[system.now]"
    28 = "[user.lastname]"
    29 = "[user.firstname]"
    30 = "[user.phone]"
    31 = "[user.mail]"
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 2).Value = $newValues[$row]
}

$ws.Range("B29").Select()
